$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5080.4
$ws.Range("I32").Value = 2400
$ws.Range("J32").Value = 5750.5
$ws.Range("K32").Value = 2400
$ws.Range("L32").Value = 5750.5
$ws.Range("M32").Value = -2074
$ws.Range("N32").Value = -6402.5

$ws.Range("H106").Value = 2046.9231
$ws.Range("I106").Value = 2101.3635
$ws.Range("K106").Value = 2101.3635
$ws.Range("M106").Value = -1470.3635

$ws.Range("H131").Value = 5166.533
$ws.Range("J131").Value = 23333
$ws.Range("L131").Value = 69999
$ws.Range("N131").Value = -80079

$ws.Range("H132").Value = 5360.788
$ws.Range("I132").Value = 1593.6923
$ws.Range("K132").Value = 4781.0769
$ws.Range("M132").Value = -2251.0769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 8065441.5
$ws.Range("I74").Value = 10000933
$ws.Range("K74").Value = 10000933
$ws.Range("M74").Value = -10000059

$ws.Range("H77").Value = 8065441.5
$ws.Range("I77").Value = 10000933
$ws.Range("K77").Value = 50004665
$ws.Range("M77").Value = -50000297

$ws.Range("H132").Value = 44398.938
$ws.Range("I132").Value = 114391
$ws.Range("J132").Value = 12584.363
$ws.Range("K132").Value = 343173
$ws.Range("L132").Value = 37753.089
$ws.Range("M132").Value = -340643
$ws.Range("N132").Value = -42813.089

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1959.85
$ws.Range("J107").Value = 2286.3845
$ws.Range("L107").Value = 2286.3845
$ws.Range("N107").Value = -6126.3845

$ws.Range("H134").Value = 2569.1333
$ws.Range("I134").Value = 1809.4615
$ws.Range("K134").Value = 5428.3845
$ws.Range("M134").Value = -2893.3845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2204.7654
$ws.Range("I31").Value = 1785.7028
$ws.Range("J31").Value = 6634.857
$ws.Range("K31").Value = 1785.7028
$ws.Range("L31").Value = 6634.857
$ws.Range("M31").Value = -1490.7028
$ws.Range("N31").Value = -7224.857

$ws.Range("H34").Value = 2204.7654
$ws.Range("I34").Value = 1785.7028
$ws.Range("J34").Value = 6634.857
$ws.Range("K34").Value = 1785.7028
$ws.Range("L34").Value = 6634.857
$ws.Range("M34").Value = -1583.7028
$ws.Range("N34").Value = -7038.857

$ws.Range("H38").Value = 19999
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

$ws.Range("H46").Value = 19999
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H86").Value = 15439.5
$ws.Range("J86").Value = 15439.5
$ws.Range("L86").Value = 15439.5
$ws.Range("N86").Value = -17685.5

$ws.Range("H89").Value = 15439.5
$ws.Range("J89").Value = 15439.5
$ws.Range("L89").Value = 77197.5
$ws.Range("N89").Value = -88429.5

$ws.Range("H132").Value = 11116911
$ws.Range("I132").Value = 15153661
$ws.Range("J132").Value = 15849.875
$ws.Range("K132").Value = 45460983
$ws.Range("L132").Value = 47549.625
$ws.Range("M132").Value = -45458453
$ws.Range("N132").Value = -52609.625

$ws.Range("H141").Value = 116551.37
$ws.Range("J141").Value = 123084.1
$ws.Range("L141").Value = 123084.1
$ws.Range("N141").Value = -133444.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 41
$ws.Range("I2").Value = 26.333334
$ws.Range("J2").Value = 52
$ws.Range("K2").Value = 158.000004
$ws.Range("L2").Value = 312
$ws.Range("M2").Value = -45.00000399999999
$ws.Range("N2").Value = -538

$ws.Range("H5").Value = 833.625
$ws.Range("I5").Value = 414.83334
$ws.Range("K5").Value = 1244.50002
$ws.Range("M5").Value = -1132.50002

$ws.Range("H135").Value = 833.625
$ws.Range("I135").Value = 414.83334
$ws.Range("K135").Value = 3733.50006
$ws.Range("M135").Value = -1198.50006

$ws.Range("H137").Value = 121430030
$ws.Range("J137").Value = 25001496
$ws.Range("L137").Value = 75004488
$ws.Range("N137").Value = -75014688

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 17248148
$ws.Range("I102").Value = 21745480
$ws.Range("K102").Value = 21745480
$ws.Range("M102").Value = -21743858

$ws.Range("H107").Value = 666.5714
$ws.Range("I107").Value = 442.6
$ws.Range("J107").Value = 791
$ws.Range("K107").Value = 442.6
$ws.Range("L107").Value = 791
$ws.Range("M107").Value = 1477.4
$ws.Range("N107").Value = -4631

$ws.Range("H122").Value = 247725.19
$ws.Range("I122").Value = 347060.03
$ws.Range("J122").Value = 7666
$ws.Range("K122").Value = 1041180.09
$ws.Range("L122").Value = 22998
$ws.Range("M122").Value = -1038730.09
$ws.Range("N122").Value = -27898

$ws.Range("H126").Value = 3865
$ws.Range("I126").Value = 2189.3333
$ws.Range("J126").Value = 7216.3335
$ws.Range("K126").Value = 6567.999899999999
$ws.Range("L126").Value = 21649.0005
$ws.Range("M126").Value = -4097.999899999999
$ws.Range("N126").Value = -26589.0005

$ws.Range("H132").Value = 1338.8823
$ws.Range("I132").Value = 787.4
$ws.Range("K132").Value = 2362.2
$ws.Range("M132").Value = 167.8000000000002

$ws.Range("H134").Value = 461968.84
$ws.Range("J134").Value = 461968.84
$ws.Range("L134").Value = 1385906.52
$ws.Range("N134").Value = -1390976.52

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 25439.75
$ws.Range("J56").Value = 34703.8
$ws.Range("L56").Value = 34703.8
$ws.Range("N56").Value = -36085.8

$ws.Range("H132").Value = 2689.37
$ws.Range("I132").Value = 2663.065
$ws.Range("J132").Value = 2777.4348
$ws.Range("K132").Value = 7989.195
$ws.Range("L132").Value = 8332.304400000001
$ws.Range("M132").Value = -5459.195
$ws.Range("N132").Value = -13392.3044

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 68997
$ws.Range("J46").Value = 68997
$ws.Range("L46").Value = 68997
$ws.Range("N46").Value = -69459

$ws.Range("H107").Value = 864.9167
$ws.Range("I107").Value = 658
$ws.Range("J107").Value = 1012.7143
$ws.Range("K107").Value = 1974
$ws.Range("L107").Value = 3038.1429
$ws.Range("M107").Value = -54
$ws.Range("N107").Value = -6878.1429

$ws.Range("H132").Value = 1877.5397
$ws.Range("I132").Value = 866.3333
$ws.Range("J132").Value = 2115.4707
$ws.Range("K132").Value = 2598.9999
$ws.Range("L132").Value = 6346.4121
$ws.Range("M132").Value = -68.9998999999998
$ws.Range("N132").Value = -11406.4121

$ws.Range("H134").Value = 68997
$ws.Range("J134").Value = 68997
$ws.Range("L134").Value = 206991
$ws.Range("N134").Value = -212061

$ws.Range("H138").Value = 98031
$ws.Range("J138").Value = 98046.5
$ws.Range("L138").Value = 98046.5
$ws.Range("N138").Value = -108326.5

$ws.Range("H140").Value = 92409
$ws.Range("J140").Value = 92409
$ws.Range("L140").Value = 92409
$ws.Range("N140").Value = -102769
